# sua chiet khau cua sale phu va update chien luoc chay tinh luong theo gio
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Ngay cong (work days)
$ws.Range("B2").Value = 21

# Phu cap (allowance)
$ws.Range("B3").Value = 735000

# Luong co ban tai LONG XUYEN (base salary)
$ws.Range("B12").Value = 3000000

# Tong luong tai LONG XUYEN (total salary at LONG XUYEN)
$ws.Range("B29").Value = 3915000

# Tong luong (grand total)
$ws.Range("B31").Value = 3915000
